$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.519.89'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.640.06'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9999'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.20'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3790'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.61'
$ws.Range("E8").Value = '  -2.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3628'
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08195'
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.237'
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9982'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.468'
$ws.Range("E14").Value = '  -3.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.406'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("D17").Value = '1.633.78'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.38'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.587'
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  -4.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9991'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  -3.11%  '
$ws.Range("D24").Value = '23.531.71'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.500'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.060'
$ws.Range("E26").Value = '  -5.77%  '
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.45'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.252'
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.30'
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("D31").Value = '1.815.64'
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.186'
$ws.Range("E32").Value = '  -5.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.668'
$ws.Range("E33").Value = '  -5.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.066'
$ws.Range("E34").Value = '  +9.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.40'
$ws.Range("E35").Value = '  +3.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02769'
$ws.Range("E36").Value = '  -4.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2501'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08777'
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07126'
$ws.Range("E39").Value = '  -2.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.040'
$ws.Range("E40").Value = '  -5.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7068'
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("E42").Value = '  -3.49%  '
$ws.Range("E43").Value = '  -4.36%  '
$ws.Range("E44").Value = '  -4.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6556'
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9989'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.289'
$ws.Range("E47").Value = '  -4.17%  '
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07980'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.80'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.193'
$ws.Range("E51").Value = '  -2.88%  '
